$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.769.63"
$ws.Range("E2").Value = "  +8.15%  "
$ws.Range("D3").Value = "1.777.58"
$ws.Range("E3").Value = "  +4.09%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'225.21"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'0.559"
$ws.Range("E6").Value = "  +4.59%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'30.72"
$ws.Range("E8").Value = "  +2.92%  "
$ws.Range("D9").Value = "'46.53"
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("E10").Value = "  +3.38%  "
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("D12").Value = "'0.0923"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "2.029.27"
$ws.Range("E13").Value = "  +4.18%  "
$ws.Range("D14").Value = "1.774.65"
$ws.Range("E14").Value = "  +3.99%  "
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "33.743.09"
$ws.Range("E16").Value = "  +8.18%  "
$ws.Range("D17").Value = "'9.99"
$ws.Range("E17").Value = "  -3.91%  "
$ws.Range("D18").Value = "'4.19"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "'68.54"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").Value = "'251.88"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'10.29"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").Value = "'4.19"
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").Value = "'158.81"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'16.50"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  +2.68%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "'3.81"
$ws.Range("E31").Value = "  +3.50%  "
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("E33").Value = "  +3.02%  "
$ws.Range("D34").Value = "'3.56"
$ws.Range("E34").Value = "  +5.04%  "
$ws.Range("D35").Value = "'1.83"
$ws.Range("E35").Value = "  +4.73%  "
$ws.Range("D36").Value = "1.484.92"
$ws.Range("D37").Value = "'1.06"
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("D38").Value = "'0.636"
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("E39").Value = "  +2.58%  "
$ws.Range("D40").Value = "'83.37"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("D42").Value = "'2.70"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "'0.885"
$ws.Range("E43").Value = "  +3.79%  "
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").Value = "'0.0513"
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("E46").Value = "  +3.27%  "
$ws.Range("D47").Value = "1.928.72"
$ws.Range("E47").Value = "  +5.19%  "
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'11.71"
$ws.Range("E50").Value = "  +12.98%  "
$ws.Range("D51").Value = "'50.74"
$ws.Range("E51").Value = "  -3.20%  "
